# The "n2" and "o1" placeholder labels in row 4 (columns B and D) are
# cleared out, but left as an explicit empty *text* entry (quote-prefixed)
# rather than a truly blank cell - matching how a user would hit
# "'" + Delete / clear the text while keeping it typed as text in Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B4").Value = "'"
$ws.Range("D4").Value = "'"
